# Rename PC_Center_Specific_Code batch names from the old
# "center_nameN" convention to the new "hms batch name N" convention,
# then mark the column as the one being reviewed/hidden from the
# detail & list views (per commit message) by selecting its populated
# cells.

$wb = $excel.ActiveWorkbook

# --- Sheet1: sample_primary_cell_batch data -------------------------------
$ws = $wb.Worksheets.Item("Sheet1")

# Column C holds PC_Center_Specific_Code (header is in row 2; data starts
# row 3). Walk the used range and rewrite any "center_nameX" value to
# "hms batch name X".
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -like "center_name*") {
        $newVal = $val -replace "^center_name", "hms batch name "
        $cell.Value = $newVal
    }
}

# Highlight/select the populated Center_Specific_Code cells (C3:C9, C11:C33
# — C10 is blank) to flag the column being hidden from the detail/list views.
$ws.Range("C3:C9,C11:C33").Select()

# --- Sheet4: secondary sheet referencing the same header layout -----------
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Activate()
$ws4.Range("A1").Select()
